$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.197644472122192
$ws.Range("B1").Value = 2.2195143699646
$ws.Range("C1").Value = 10.61748123168945
$ws.Range("D1").Value = 2.57518196105957
$ws.Range("E1").Value = 1.22445547580719
